# Build site at 2022-01-09 00:29:46 UTC
# Adds a "Docentes responsáveis:" row with three lecturer names to the
# LOM3104 course sheet, inserted right after the "Objectives:" row and
# before the "Programa resumido:" row. Everything below shifts down by
# four rows (old row 12 -> new row 16, ... old row 22 -> new row 26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 12 (old row 12 "Programa resumido:"
# and everything after it shifts down to make room).
$ws.Rows("12:15").Insert()

# The insert stamps column A's sheet-level column style onto the new
# blank rows (A13:A15) even though no data belongs there. Clear those
# cells completely so they drop out of the saved sheet, matching rows
# that never had an "A" cell (e.g. row 13-15 only have B/C content).
$ws.Range("A13:A15").Clear()

# Row 12: only column A, the new section header.
$ws.Range("A12").Value = "Docentes responsáveis:"

# Rows 13-15 carry the three lecturer names in both B (plain) and C (the
# "modified data" column). Copy the column B/C formatting from the row
# that now sits at 16 (the old row 12, "Programa resumido:") so the new
# cells pick up the existing style indices instead of minting new ones.
$ws.Range("B16:C16").Copy()
$ws.Range("B13:C15").PasteSpecial(-4122)

$ws.Range("B13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"
$ws.Range("C13").Value = "7459752 - Maria Ismenia Sodero Toledo Faria"

$ws.Range("B14").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C14").Value = "2166002 - Sandra Giacomin Schneider"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"
